$wb = $excel.ActiveWorkbook

# --- Rename Sheet1 -> linkspythoncode ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Name = "linkspythoncode"

# --- linkspythoncode sheet data ---
$ws1.Range("A1").Value = "links"
$ws1.Range("B1").Value = "pythonCode"
$ws1.Range("C1").Value = "Result"
$ws1.Range("A2").Value = "arrays-in-python"
$ws1.Range("B2").Value = "print('hello')"
$ws1.Range("C2").Value = "hello"
$ws1.Range("A3").Value = "arrays-using-list"
$ws1.Range("B3").Value = "print('hello')"
$ws1.Range("C3").Value = "hello"
$ws1.Range("A4").Value = "basic-operations-in-lists"
$ws1.Range("B4").Value = "print('hello')"
$ws1.Range("C4").Value = "hello"
$ws1.Range("A5").Value = "applications-of-array"
$ws1.Range("B5").Value = "print('hello')"
$ws1.Range("C5").Value = "hello"
$ws1.Range("A6").NumberFormat = "General"
$ws1.Range("B6").NumberFormat = "General"
$ws1.Range("C6").NumberFormat = "General"
$ws1.Range("A7").Value = "arrays-in-python"
$ws1.Range("B7").Value = "hello"
$ws1.Range("C7").Value = "NameError: name 'hello' is not defined on line 1"
$ws1.Range("A8").Value = "arrays-using-list"
$ws1.Range("B8").Value = "hello"
$ws1.Range("C8").Value = "NameError: name 'hello' is not defined on line 1"
$ws1.Range("A9").Value = "basic-operations-in-lists"
$ws1.Range("B9").Value = "hello"
$ws1.Range("C9").Value = "NameError: name 'hello' is not defined on line 1"
$ws1.Range("A10").Value = "applications-of-array"
$ws1.Range("B10").Value = "hello"
$ws1.Range("C10").Value = "NameError: name 'hello' is not defined on line 1"

# --- Column widths on linkspythoncode (best achievable approximation) ---
$ws1.Columns.Item(1).ColumnWidth = 20.5
$ws1.Columns.Item(2).ColumnWidth = 10.6666666667

# --- pythonCode sheet data ---
$ws2 = $wb.Worksheets.Item("pythonCode")
$ws2.Range("A1").Value = "pythonCode"
$ws2.Range("B1").Value = "Result"
$ws2.Range("A2").Value = "print(`"hello`");"
$ws2.Range("B2").Value = "hello"
$ws2.Range("A3").Value = "hello"
$ws2.Range("B3").Value = "hello"
$ws2.Range("A4").Value = "def search(input_list, num):`n  if(num in input_list):`n    print(`"Element Found`")`n  else:`n    print(`"Not Found`")`nsearch([12, 23, 45, 67, 6, 90] , 12)"
$ws2.Range("B4").Value = "Element Found"
$ws2.Range("A5").Value = "def search(input_list, num):`n  if(num in input_list):`n    print(`"Element Found`")`n  else:`n    print(`"Not Found`")`nsearch([12, 23, 45, 67, 6, 90] , 12)"
$ws2.Range("B5").Value = "submission success"
$ws2.Range("A6").Value = "def findMaxConsecutiveOnes(nums) :`n  result = 0`n  count = 0`n  for i in nums:`n    if i == 0:`n      count = 0 `n    else:`n      count+= 1`n      result = max(result, count) `n    `n  return result`n `t`nprint(findMaxConsecutiveOnes([1,0,1,1,0,1]))"
$ws2.Range("B6").NumberFormat = "@"
$ws2.Range("B6").Value = "2"
$ws2.Range("A7").Value = "def findMaxConsecutiveOnes(nums) :`n  result = 0`n  count = 0`n  for i in nums:`n    if i == 0:`n      count = 0 `n    else:`n      count+= 1`n      result = max(result, count) `n    `n  return result`n `t`nprint(findMaxConsecutiveOnes([1,0,1,1,0,1]))"
$ws2.Range("B7").Value = "submission success"
$ws2.Range("A8").Value = "def findNumbers(nums):`n  c=0`n  for i in nums:`n    j=str(i)`n    x=len(j)`n    if x%2==0:`n       c=c+1`n  print c`n  return c`nfindNumbers([555,901,482,1771])"
$ws2.Range("B8").NumberFormat = "@"
$ws2.Range("B8").Value = "1"
$ws2.Range("A9").Value = "def findNumbers(nums):`n  c=0`n  for i in nums:`n    j=str(i)`n    x=len(j)`n    if x%2==0:`n       c=c+1`n  print c`n  return c`nfindNumbers([555,901,482,1771])"
$ws2.Range("B9").Value = "submission success"
$ws2.Range("A10").Value = "def sortedSquares(nums):`n  squares_list = []`n  for i in range(0, len(nums)):`n    square = nums[i] * nums[i];`n    squares_list.append(square)`n  sorted_squares_list = sorted(squares_list)`n  print sorted_squares_list;`n  return sorted_squares_list;`nsortedSquares([-7,-3,2,3,11])"
$ws2.Range("B10").Value = "[4, 9, 9, 49, 121]"
$ws2.Range("A11").Value = "def sortedSquares(nums):`n  squares_list = []`n  for i in range(0, len(nums)):`n    square = nums[i] * nums[i];`n    squares_list.append(square)`n  sorted_squares_list = sorted(squares_list)`n  print sorted_squares_list;`n  return sorted_squares_list;`nsortedSquares([-7,-3,2,3,11])"
$ws2.Range("B11").Value = "submission success"
$ws2.Range("A12").Value = "def search(input_list, num):`n  if(num in input_list):`n    print(`"Element Found`n  else:`n    print(`"Not Found`")`nsearch([12, 23, 45, 67, 6, 90] , 12)"
$ws2.Range("B12").Value = "SyntaxError: bad input on line 3"
$ws2.Range("A13").Value = "def findMaxConsecutiveOnes(nums) :`n  result = 0`n  count = 0`n for i in nums:`n    if i == 0:`n      count = 0 `n    else:`n      count+= 1`n      result = max(result, count) `n  return result`nprint(findMaxConsecutiveOnes([1,0,1,1,0,1]))"
$ws2.Range("B13").Value = "SyntaxError: unindent does not match any outer indentation level on line 4"
$ws2.Range("A14").Value = "def findNumbers(nums):`n  c=0`n  for i in nums:`n    j=str(i)`n    x=len(j)`n    if x%2==0:`n       c=c+1`n  print c`n return c`nfindNumbers([555,901,482,1771]"
$ws2.Range("B14").Value = "SyntaxError: unindent does not match any outer indentation level on line 9"
$ws2.Range("A15").Value = "def sortedSquares(nums)`n  squares_list = []`n  for i in range(0, len(nums)):`n    square = nums[i] * nums[i];`n    squares_list.append(square)`n  sorted_squares_list = sorted(squares_list)`n  print sorted_squares_list;`n  return sorted_squares_list;`nsortedSquares([-7,-3,2,3,11])"
$ws2.Range("B15").Value = "SyntaxError: bad input on line 1"
$ws2.Range("B16").Value = "Enter Code before Run"
$ws2.Range("B17").Value = "No tests were collected"
